$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.909.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.61%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.811.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.60%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'310.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.93%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  -0.01%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.4631"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +3.39%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3757"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.06%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.07416"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.30%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.8785"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.77%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'20.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.61%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.815.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.40%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'5.360"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.79%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'6.550"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.08%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.07051"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.82%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'91.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.55%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.03%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'  -0.63%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.03%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'14.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.05%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'26.907.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.65%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'5.313"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.06%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'10.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.48%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'2.026.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.41%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'1.924"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.14%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  +0.02%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'18.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.50%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'2.155"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -9.15%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'5.308"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.04%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'116.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.88%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.08908"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.75%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.7715"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.65%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.167"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.56%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'4.488"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.48%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'2.900"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.75%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  -0.01%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  +0.53%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.01959"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.87%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'2.437"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +4.57%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.05239"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.69%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'7.275"
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").Value = "'0.5359"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.75%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'2.908"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.67%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.1663"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.63%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'8.580"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.86%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.5078"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.80%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'10.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.90%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'104.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.41%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  -0.01%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'1.667"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.03%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.06330"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.80%  "
$ws.Range("E51").Style = "Normal"
